# Append a new data row (row 3) below the existing firstname/lastname table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Alkesh "
$ws.Range("B3").Value = "Rathore"

# Match the saved workbook view state: selection moves to the newly entered cell.
$null = $ws.Range("B3").Select()
